$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated cryptos list (price + 1h volume refresh from GitHub Actions run) ---
# Price column (D) holds plain text values (e.g. "505.50"), so force the
# number format to Text before writing, then restore the default "Normal"
# style afterwards so we don't leave a stray text-format style behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.927.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.291.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.43%  "

$ws.Range("E7").Value = "  -0.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.312.75"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.79%  "

$ws.Range("E10").Value = "  +2.13%  "

$ws.Range("E11").Value = "  +1.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.96%  "

$ws.Range("E13").Value = "  +0.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.700.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "54.922.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.79%  "

$ws.Range("E17").Value = "  +1.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.301.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.46%  "

$ws.Range("E19").Value = "  +2.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.36%  "

# --- Rows 21 & 22: Uniswap and BitcoinCash swapped rank positions ---
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.56%  "

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "310.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.77%  "

$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("E27").Value = "  +2.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.63%  "

$ws.Range("E29").Value = "  +3.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0707"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.83%  "

$ws.Range("E32").Value = "  +4.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.05"
$ws.Range("D34").Style = "Normal"

$ws.Range("E35").Value = "  -0.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.924"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.49%  "

$ws.Range("E37").Value = "  +2.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("E40").Value = "  +2.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "134.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.82%  "

$ws.Range("E43").Value = "  +1.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "261.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.86%  "

$ws.Range("E46").Value = "  +1.98%  "

$ws.Range("E47").Value = "  +1.75%  "

$ws.Range("E48").Value = "  +0.24%  "

$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0211"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.98%  "
